$d = $word.ActiveDocument

function New-WordParagraphXml($innerRunsXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerRunsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. "color" -> "colour" in the zombie-game description paragraph, and drop
#    the spell-check proofErr wrapper that surrounded the misspelling.
# ---------------------------------------------------------------------------
$colorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*different color circle*") {
        $colorPara = $p
        break
    }
}

if ($colorPara -ne $null) {
    $runsXml = '<w:r><w:t xml:space="preserve">It is essentially a zombie infection game where humans and zombies are represented with a different </w:t></w:r>' +
               '<w:r><w:t>colour</w:t></w:r>' +
               '<w:r><w:t xml:space="preserve"> circle</w:t></w:r>' +
               '<w:r><w:t>. This can be changed later if needed.</w:t></w:r>'
    $colorPara.Range.InsertXML((New-WordParagraphXml $runsXml))
}

# ---------------------------------------------------------------------------
# 2. After the "Barricades can be destroyed..." paragraph, add a blank
#    paragraph followed by a new paragraph describing the "humans that can
#    shoot" idea.
# ---------------------------------------------------------------------------
$barricadePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Barricades can be destroyed by zombies depending on the health stat of the barricade.*") {
        $barricadePara = $p
        break
    }
}

if ($barricadePara -ne $null) {
    $barricadePara.Range.InsertParagraphAfter()
    $blankPara = $barricadePara.Next()

    $blankPara.Range.InsertParagraphAfter()
    $newPara = $blankPara.Next()
    $newPara.Style = $blankPara.Style

    $runsXml = '<w:r><w:t xml:space="preserve">If time is left, then </w:t></w:r>' +
               '<w:r><w:t>humans that</w:t></w:r>' +
               '<w:r><w:t xml:space="preserve"> can shoot will probably be implemented.</w:t></w:r>' +
               '<w:r><w:t xml:space="preserve"> This will give the humans a chance to </w:t></w:r>' +
               '<w:r><w:t>survive</w:t></w:r>' +
               '<w:r><w:t>.</w:t></w:r>'
    $newPara.Range.InsertXML((New-WordParagraphXml $runsXml))
}
